$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format Price/Volume columns as Text so numeric-looking strings
# (e.g. "420.24", "0.0000369") are stored as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '66.526.79'
$ws.Range("E2").Value = '  +4.59%  '
$ws.Range("D3").Value = '3.677.50'
$ws.Range("E3").Value = '  +5.68%  '
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").Value = '420.24'
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("D6").Value = '128.79'
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("D7").Value = '3.678.93'
$ws.Range("E7").Value = '  +5.96%  '
$ws.Range("D8").Value = '0.634'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = '0.752'
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("D11").Value = '0.178'
$ws.Range("E11").Value = '  +11.28%  '
$ws.Range("D12").Value = '0.0000369'
$ws.Range("E12").Value = '  +57.85%  '
$ws.Range("D13").Value = '41.77'
$ws.Range("E13").Value = '  -1.57%  '
$ws.Range("D14").Value = '9.79'
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").Value = '4.244.93'
$ws.Range("E15").Value = '  +4.97%  '
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '20.40'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = '3.685.51'
$ws.Range("E18").Value = '  +5.65%  '
$ws.Range("D19").Value = '12.97'
$ws.Range("E19").Value = '  +4.50%  '
$ws.Range("D20").Value = '1.11'
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").Value = '66.484.95'
$ws.Range("E21").Value = '  +4.63%  '
$ws.Range("D22").Value = '443.87'
$ws.Range("E22").Value = '  -3.11%  '
$ws.Range("D23").Value = '15.37'
$ws.Range("E23").Value = '  +15.76%  '
$ws.Range("D24").Value = '87.71'
$ws.Range("E24").Value = '  -2.84%  '
$ws.Range("D25").Value = '3.10'
$ws.Range("E25").Value = '  -5.70%  '
$ws.Range("D26").Value = '37.33'
$ws.Range("E26").Value = '  +11.30%  '
$ws.Range("D27").Value = '10.22'
$ws.Range("E27").Value = '  +0.95%  '
$ws.Range("D28").Value = '3.26'
$ws.Range("E28").Value = '  -1.23%  '
$ws.Range("D29").Value = '4.96'
$ws.Range("E29").Value = '  +4.28%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '2.76'
$ws.Range("E30").Value = '  +2.51%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '12.32'
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").Value = '0.118'
$ws.Range("E32").Value = '  +4.53%  '
$ws.Range("D33").Value = '7.14'
$ws.Range("E33").Value = '  -4.92%  '
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").Value = '41.75'
$ws.Range("E34").Value = '  +4.26%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '0.160'
$ws.Range("E35").Value = '  -4.16%  '
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").Value = '56.17'
$ws.Range("E37").Value = '  -3.01%  '
$ws.Range("D38").Value = '0.0481'
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("B39").Value = 'ThetaToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D39").Value = '3.00'
$ws.Range("E39").Value = '  +28.57%  '
$ws.Range("D40").Value = '0.0₃0697'
$ws.Range("E40").Value = '  +7.18%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.144'
$ws.Range("E41").Value = '  +4.46%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '28.22'
$ws.Range("E42").Value = '  +29.19%  '
$ws.Range("D43").Value = '0.995'
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("D44").Value = '3.40'
$ws.Range("E44").Value = '  +1.79%  '
$ws.Range("D45").Value = '144.42'
$ws.Range("E45").Value = '  -1.03%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '2.66'
$ws.Range("E46").Value = '  -6.39%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '2.08'
$ws.Range("E47").Value = '  +4.05%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '2.88'
$ws.Range("E48").Value = '  -7.00%  '
$ws.Range("D49").Value = '4.24'
$ws.Range("E49").Value = '  -5.27%  '
$ws.Range("D50").Value = '0.303'
$ws.Range("E50").Value = '  -4.53%  '
$ws.Range("D51").Value = '0.156'
$ws.Range("E51").Value = '  +12.21%  '

# Restore original (unstyled) formatting for the Price/Volume columns.
$ws.Range("D2:E51").ClearFormats()
